$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the year value from 2019 to 2015 in cell A2
$ws.Range("A2").Value = 2015

# Update the active selection to A3 (matches saved sheetView selection)
$ws.Range("A3").Select()
